$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.146.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5041"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3921"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09683"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.139"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.499"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.867.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.001"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.410"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001127"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06614"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.155"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.210.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.282"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.534"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.081.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1059"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.064"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.622"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.624"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.558"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06731"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02385"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6350"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.974"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.178"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6012"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.657"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.260"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.990"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("E51").Value = "  +0.66%  "
